$wb = $excel.ActiveWorkbook
$srcSheet = $wb.Worksheets.Item("2025-07-14")
$srcSheet.Copy([System.Reflection.Missing]::Value, $srcSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2025-07-15"

$rows = @(
    @(1, "帰ってください！ 阿久津さん", "長岡太一(著者)", "第192話"),
    @(2, "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～", "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり", "第５０話　雌雄を決する器用貧乏（２）"),
    @(3, "【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！", "島知宏 音速炒飯 有都あらゆる", "第２２食　ユクシーさんの覚悟、すごいのですわ！（２）"),
    @(4, "魔導具師ダリヤはうつむかない ～Dahliya Wilts No More～", "漫画：住川惠 原作：甘岸久弥(｢魔導具師ダリヤはうつむかない ～今日から自由な職人ライフ～｣MFブックス刊) キャラクター原案：景、駒田ハチ", "第47話 魔導具師とつながれたもの②"),
    @(5, "十年目、帰還を諦めた転移者はいまさら主人公になる", "原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう", "第１８話②"),
    @(6, "異世界でも無難に生きたい症候群", "原作：安泰（一二三書房刊） 漫画：笹峰コウ キャラクター原案：ひたきゆう", "第30話①"),
    @(7, "ギャルとダンジョンと周回遅れの探索英雄譚", "漫画家： 水田ケンジ 原作：榊一郎 キャラクター原案：黒獅子", "第1話"),
    @(8, "賢者の孫", "緒方俊輔(漫画) 吉岡剛(原作) 菊池政治(キャラクター原案)", "第94話-2"),
    @(9, "めっちゃ召喚された件 THE COMIC", "漫画：六甲島カモメ 原作：さいとうさ キャラクター原案：ツグトク", "第46話②"),
    @(10, "レベル１だけどユニークスキルで最強です", "漫画：真綿 原作：三木なずな キャラクター原案：すばち", "第７２話　インドール救援作戦 大成功!?（１）"),
    @(11, "異世界建国記", "ＫＯＩＺＵＭＩ(漫画) 桜木桜(原作) 屡那(キャラクター原案)", "第77話"),
    @(12, "外れスキル『レベルアップ』のせいでパーティーを追放された少年は、レベルを上げて物理で殴る", "しんこせい 大橋ウルオ てんまそ", "第18話 ミヤダケ（前編）"),
    @(13, "大学入学時から噂されていた美少女三姉妹、生き別れていた義妹だった。", "こすずめ(著者) 夏乃実(原作) ポメ(キャラクター原案)", "今週はスペシャルイラスト！！"),
    @(14, "無能と呼ばれた『精霊たらし』～実は異能で、精霊界では伝説的ヒーローでした～＠COMIC", "原作：佐藤謙羊 漫画：タバタグランドキャニオン", "第27話「求職は総合人材ギルドで」①"),
    @(15, "人外姫様、始めました　-Free Life Fantasy Online-", "園原アオ 割田コマ 子日あきすず Sherry", "第５９話　再チャレンジと思ったら、まさかの冥府へ！（２）"),
    @(16, "器用貧乏、城を建てる～開拓学園の劣等生なのに、 上級職のスキルと魔法がすべて使えます～＠COMIC", "原作：佐藤謙羊 漫画：スガン", "第21話②「ルールの代償」"),
    @(17, "尾守つみきと奇日常。", "森下みゆ", "第58話 雷漸会長と雷。"),
    @(18, "転生少女はまず一歩からはじめたい～魔物がいるとか聞いてない！～", "原作：カヤ 漫画：岡村アユム キャラクター原案：那流", "第37話 お屋敷での日々①"),
    @(19, "魔眼の悪役に転生したので推しキャラを見守るモブを目指します", "在間りしん(漫画) 瀧岡くるじ(原作) 福きつね(キャラクター原案)", "第10話②"),
    @(20, "没落予定の貴族だけど、暇だったから魔法を極めてみた@COMIC", "漫画：秋咲りお 原作：三木なずな キャラクター原案：かぼちゃ", "第52話"),
    @(21, "レベル１から始まる召喚無双", "漫画：七桃りお 原作：白石新 キャラクター原案：夕薙", "第32話後半"),
    @(22, "狂戦士なモブ、無自覚に本編を破壊する", "漫画：佐藤良亮 原作：なるのるな キャラクター原案：霜月えいと", "第11話 ②"),
    @(23, "超!!! 天才発明令嬢のパワフル領地改革", "漫画：katoson 原作：守雨", "第46話"),
    @(24, "元最強の剣士は、異世界魔法に憧れる", "漫画：天乃ちはる 原作：紅月シン キャラクター原案：necömi", "第72話後半"),
    @(25, "アラフォーになった最強の英雄たち、再び戦場で無双する!!", "漫画：戸田泰成 原作：岸馬きらく 構成協力：森小太郎 キャラクター原案：peroshi", "第23話 後編"),
    @(26, "女神から『孵化』のスキルを授かった俺が、なぜか幻獣や神獣を従える最強テイマーになるまで", "漫画家：春夏冬 唯人 原作：まるせい", "第15話"),
    @(27, "S級パーティーから追放された狩人、実は世界最強 ～射程9999の男、帝国の狙撃手として無双する～", "漫画：カズミヤアキラ 原作：茨木野 キャラクター原案：へいろー", "第9話 ①"),
    @(28, "異世界のんびり開拓記  -平凡サラリーマン、万能自在のビルド&クラフトスキルで、気ままなスローライフ 開拓始めます! -", "漫画：しょうじひでまさ 原作：タライ和治 キャラクター原案：イシバシヨウスケ", "第22話 後編"),
    @(29, "寿司ガキ", "ichika(著者)", "第9話"),
    @(30, "不純な彼女達は懺悔しない", "ポロロッカ(著者)", "第28話"),
    @(31, "クセ強彼女は床にいざなう", "須河篤志(著者)", "第13話後半"),
    @(32, "忌み子と呼ばれた召喚士@COMIC", "漫画：コイシ 原作：緑黄色野菜 キャラクター原案：こよいみつき", "第40話 ②"),
    @(33, "ダンジョンバンド", "新挑限(著者)", "♯７SPEED TUBE ①"),
    @(34, "ぽんドロイド！ はまさん", "はれやまはれぞう(著者)", "第4話"),
    @(35, "ヘテロゲニア リンギスティコ ～異種族言語学入門～", "瀬野反人(著者)", "第56話"),
    @(36, "ちはるくんは女装をしたくない！", "翁丸ジョン", "【第20話】男装宗と交流したくない！その二"),
    @(37, "千年英雄", "原作/福島航平 作画/中村ゆきひろ", "18話②"),
    @(38, "おとなりのダウナーさんは無理させない", "瑠璃いろ(著者)", "第12話"),
    @(39, "スライムテイマーの異世界ライフ", "作画：チイチ 原作：穴の空いた靴下", "第7話(1)"),
    @(40, "魔術師クノンは見えている", "La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)", "第38話②"),
    @(41, "リビルドワールド", "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)", "第71話①"),
    @(42, "地味子な三葉さんが僕を誘惑する", "はぶらえる(著者)", "第10話前半"),
    @(43, "勇者に全部奪われた俺は勇者の母親とパーティを組みました！", "久遠まこと(著者) 石のやっさん(原作)", "第28話"),
    @(44, "サーシャちゃんとクラスメイトオタクくん", "はぐはぐ(著者)", "第5話"),
    @(45, "勇者ディンは2.5回殺せ", "ナゲク 福留のりひろ", "第4話「その死は必然」②"),
    @(46, "女友達は頼めば意外とヤらせてくれる", "ろくろ(漫画) 鏡遊(原作)", "特別イラスト"),
    @(47, "わたしのために脱ぎなさいっ！", "九郎(著者)", "第83話"),
    @(48, "淫獄団地", "搾精研究所(原作) 丈山雄為(漫画)", "第48話（後編）"),
    @(49, "悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる", "こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)", "「名もなき英雄譚」前半"),
    @(50, "冒険者ギルドが十二歳からしか入れなかったので、サバよみました。", "GUNP（漫画） KAME （原作） ox （キャラクター原案）", "第12話後半"),
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
    $ws.Cells.Item($r, 4).Value = $rows[$i][3]
}

[void]$ws.Range("A1").Select()
[void]$wb.Worksheets.Item("Sheet1").Activate()
